# Update cryptocurrency price/volume data per upstream refresh (GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.167.74"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.389.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'548.08"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'141.46"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -10.88%  "
$ws.Range("D9").Value = "2.387.33"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'5.28"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "'0.347"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "'25.42"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "2.821.54"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "'0.0000165"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "61.130.71"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "2.389.43"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("D19").Value = "'10.79"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'318.46"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'6.70"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'1.91"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'63.52"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'8.20"
$ws.Range("E26").Value = "  +5.27%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "2.509.00"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "0.0₃0928"
$ws.Range("E29").Value = "  -6.88%  "
$ws.Range("D30").Value = "'525.07"
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").Value = "'8.08"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'5.51"
$ws.Range("E37").Value = "  -6.33%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.69"
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("D39").Value = "'0.376"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").Value = "'1.83"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("D41").Value = "'18.09"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").Value = "'138.39"
$ws.Range("E42").Value = "  -6.00%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'40.30"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "'2.15"
$ws.Range("E45").Value = "  -9.24%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'140.71"
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.61"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'20.11"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("D49").Value = "'0.0519"
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("D50").Value = "'0.575"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("D51").Value = "'0.0226"
$ws.Range("E51").Value = "  -1.17%  "
